# Update the marksheet figures on the "quiz" sheet:
#  - B11 (Marking / Right count): 3  -> 5
#  - B12 (Total / Right count):   75 -> 125
#  - E12 (Total / Max "corr/total"): "74/84" -> "125/140"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 125
$ws.Range("E12").Value = "125/140"
